$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 8.785400000000001
$ws.Range("D2").Value = -7.110899999999996
$ws.Range("A3").Value = -21.85570000000001
$ws.Range("C3").Value = -10.8508
$ws.Range("D6").Value = -7.879799999999992
$ws.Range("E8").Value = 16.0795
$ws.Range("E9").Value = 17.34280000000001
$ws.Range("C12").Value = -12.4767
$ws.Range("A14").Value = -21.52760000000001
$ws.Range("A16").Value = -21.70939999999999
$ws.Range("B18").Value = 6.956099999999995
$ws.Range("D19").Value = -9.093599999999995
$ws.Range("A21").Value = -20.24689999999998
$ws.Range("A23").Value = -20.71289999999997
$ws.Range("E23").Value = 16.2366
$ws.Range("B24").Value = 6.9358
$ws.Range("C24").Value = -12.74869999999999
$ws.Range("D24").Value = -7.930599999999997
$ws.Range("A25").Value = -21.61679999999999
$ws.Range("B25").Value = 5.359899999999993
$ws.Range("C25").Value = -13.887
$ws.Range("A26").Value = -21.21399999999996
$ws.Range("E26").Value = 16.24099999999999
$ws.Range("B27").Value = 5.603299999999998
$ws.Range("D27").Value = -8.732900000000004
$ws.Range("A29").Value = -20.81409999999997
$ws.Range("B30").Value = 5.772499999999999
$ws.Range("D30").Value = -7.633800000000003
$ws.Range("B31").Value = 5.044500000000002
$ws.Range("D31").Value = -8.723300000000002
$ws.Range("D33").Value = -7.701799999999997
$ws.Range("E37").Value = 16.53280000000001
$ws.Range("B39").Value = 9.540900000000001
$ws.Range("A40").Value = -20.228
$ws.Range("C41").Value = -12.6019
$ws.Range("B42").Value = 10.07549999999999
$ws.Range("D42").Value = -9.143299999999995
$ws.Range("B48").Value = 5.173600000000001
$ws.Range("E48").Value = 17.7938
$ws.Range("C50").Value = -13.35199999999998
$ws.Range("B51").Value = 5.500300000000001
$ws.Range("B52").Value = 5.213700000000001
$ws.Range("A53").Value = -21.5025
$ws.Range("C53").Value = -10.4472
$ws.Range("E54").Value = 16.61060000000001
$ws.Range("B55").Value = 6.662699999999994
$ws.Range("D55").Value = -8.0495
$ws.Range("B56").Value = 5.703000000000002
$ws.Range("C56").Value = -12.19579999999999
$ws.Range("A57").Value = -21.9813
$ws.Range("B57").Value = 5.565999999999998
$ws.Range("C57").Value = -12.5825
$ws.Range("C58").Value = -13.17529999999999
$ws.Range("D58").Value = -8.363699999999994
$ws.Range("A59").Value = -22.31350000000001
$ws.Range("B60").Value = 5.470499999999998
$ws.Range("C61").Value = -12.97659999999999
$ws.Range("E62").Value = 16.5914
$ws.Range("C63").Value = -11.727
$ws.Range("C64").Value = -12.21799999999999
$ws.Range("A65").Value = -21.92739999999998
$ws.Range("D65").Value = -8.0374
$ws.Range("E65").Value = 17.00119999999998
$ws.Range("E66").Value = 17.6769
$ws.Range("A69").Value = -21.55250000000001
$ws.Range("C70").Value = -10.97530000000001
$ws.Range("D70").Value = -8.170400000000006
$ws.Range("C72").Value = -11.67810000000001
$ws.Range("B73").Value = 8.465299999999997
$ws.Range("B74").Value = 9.2805
$ws.Range("D74").Value = -8.689700000000002
$ws.Range("D75").Value = -8.0809
$ws.Range("E75").Value = 16.6288
$ws.Range("A79").Value = -20.45100000000001
$ws.Range("E81").Value = 16.88629999999999
$ws.Range("A83").Value = -22.00269999999999
$ws.Range("D83").Value = -8.543599999999993
$ws.Range("D84").Value = -8.809899999999999
$ws.Range("C86").Value = -13.00469999999999
$ws.Range("D86").Value = -7.923699999999992
$ws.Range("B89").Value = 5.031900000000002
$ws.Range("C89").Value = -11.4453
$ws.Range("E89").Value = 17.21790000000001
$ws.Range("B90").Value = 5.560599999999997
$ws.Range("A91").Value = -21.2647
$ws.Range("B92").Value = 5.375999999999992
$ws.Range("A93").Value = -20.98019999999997
$ws.Range("E94").Value = 19.12850000000001
$ws.Range("D96").Value = -8.497299999999997
$ws.Range("D97").Value = -8.022599999999997
$ws.Range("C98").Value = -12.03359999999999
$ws.Range("A100").Value = -21.6922
$ws.Range("C100").Value = -13.62699999999999
$ws.Range("C102").Value = -12.75180000000001
